$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Assess Crime across years in DC area" -> "Crime across years "
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(29)
$r1 = $p1.Range
$r1.Find.Execute("Assess ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

$p1 = $d.Paragraphs(29)
$r1 = $p1.Range
$r1.Find.Execute("in DC area", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------------
# Change 2: "Which district has the highest level of Crime & what type of Crime"
#           -> "High Crime or unsafe area"
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(35)
$r2 = $p2.Range
$r2.Find.Execute("Which district has the highest level of Crime & what type of Crime", $false, $false, $false, $false, $false, $true, 1, $false, "High Crime or unsafe area", 2)

# ---------------------------------------------------------------------------
# Change 3: rewrite the Homoside question and relocate the _GoBack bookmark
#           to the end of the paragraph (after all of its text).
# ---------------------------------------------------------------------------

# Drop the existing hidden "_GoBack" bookmark; it will be re-added once the
# text has been rearranged.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$p3 = $d.Paragraphs(38)
$r3 = $p3.Range
$r3.Find.Execute("What is the rate of homoside across DC, drill into gun related homoside", $false, $false, $false, $false, $false, $true, 1, $false, "Homoside across DC", 2)

$p3 = $d.Paragraphs(38)
$r3 = $p3.Range
$r3.Find.Execute("– gmap API to used to ", $false, $false, $false, $false, $false, $true, 1, $false, " – gmap API to used to ", 2)

# Temporarily mark the very end of the paragraph with a unique marker string so
# the bookmark position can be located via Find (a raw End-1 offset right
# against the paragraph mark trips an off-by-one in this host's Bookmarks.Add).
$p3 = $d.Paragraphs(38)
$endRange = $d.Range($p3.Range.End - 1, $p3.Range.End - 1)
$endRange.InsertAfter("@@MARK@@")

$p3 = $d.Paragraphs(38)
$r3 = $p3.Range
$r3.Find.Execute("@@MARK@@")
$markStart = $d.Range($r3.Start, $r3.Start)
$d.Bookmarks.Add("_GoBack", $markStart)

$p3 = $d.Paragraphs(38)
$r3 = $p3.Range
$r3.Find.Execute("@@MARK@@")
$r3.Delete()

Write-Output "done"
